# Add derived statistics columns/rows to the "Sheet2" (Districts, sheet3.xml)
# and "Sheet3" (Locality, sheet2.xml) worksheets.
#
# Workbook tab order is Sheet1, Sheet3, Sheet2. We finish by touching
# "Sheet3" last so it remains the active tab/selection, matching the
# original file's state. Shared-string writes are ordered to match the
# first-use order of each new label.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Sheet2" worksheet (districts table, rows 2-32) -> xl/worksheets/sheet3.xml
# ---------------------------------------------------------------------
$districts = $wb.Worksheets.Item("Sheet2")

# New header (column H) + per-row average for every data row.
$districts.Range("H1").Value = "AVERAGE"
for ($r = 2; $r -le 32; $r++) {
    $districts.Range("H$r").Formula = "=AVERAGE(B$r`:G$r)"
}

# Row 35: column maxima.
$districts.Range("A35").Value = "max_val_col"
$districts.Range("B35").Formula = "=MAX(B2:B32)"
$districts.Range("C35").Formula = "=MAX(C2:C32)"
$districts.Range("D35").Formula = "=MAX(D2:D32)"
$districts.Range("E35").Formula = "=MAX(E2:E32)"
$districts.Range("F35").Formula = "=MAX(F2:F32)"
$districts.Range("G35").Formula = "=MAX(G2:G32)"

# Row 36: column minima.
$districts.Range("A36").Value = "min_val_col"
$districts.Range("B36").Formula = "=MIN(B2:B32)"
$districts.Range("C36").Formula = "=MIN(C2:C32)"
$districts.Range("D36").Formula = "=MIN(D2:D32)"
$districts.Range("E36").Formula = "=MIN(E2:E32)"
$districts.Range("F36").Formula = "=MIN(F2:F32)"
$districts.Range("G36").Formula = "=MIN(G2:G32)"

# Row 37: column range (max - min). Only B37 carries the live formula;
# C37:G37 were typed in as the resulting values.
$districts.Range("A37").Value = "range_val_col"
$districts.Range("B37").Formula = "=MAX(B2:B32)-MIN(B2:B32)"
$districts.Range("C37").Value = 3.1
$districts.Range("D37").Value = 10.6
$districts.Range("E37").Value = 6.5
$districts.Range("F37").Value = 5
$districts.Range("G37").Value = 16.2

# Row 38: column averages.
$districts.Range("A38").Value = "AVERAGE_COL_VAL"
$districts.Range("B38").Formula = "=AVERAGE(B2:B32)"
$districts.Range("C38").Formula = "=AVERAGE(C2:C32)"
$districts.Range("D38").Formula = "=AVERAGE(D2:D32)"
$districts.Range("E38").Formula = "=AVERAGE(E2:E32)"
$districts.Range("F38").Formula = "=AVERAGE(F2:F32)"
$districts.Range("G38").Formula = "=AVERAGE(G2:G32)"
$districts.Range("H38").Formula = "=AVERAGE(H2:H32)"

# Combined-average headers + values.
$districts.Range("I1").Value = "AVERAGE_WOMEN"
$districts.Range("J1").Value = "AVERAGE_MEN"
$districts.Range("I38").Formula = "=AVERAGE(B38:D38)"
$districts.Range("J38").Formula = "=AVERAGE(E38:G38)"

$districts.Activate() | Out-Null
$districts.Range("G1").Select() | Out-Null

# ---------------------------------------------------------------------
# "Sheet3" worksheet (locality table, rows 2-3) -> xl/worksheets/sheet2.xml
# ---------------------------------------------------------------------
$locality = $wb.Worksheets.Item("Sheet3")

$locality.Range("H1").Value = "Avg_total"
$locality.Range("H2").Formula = "=AVERAGE(B2:G2)"
$locality.Range("H3").Formula = "=AVERAGE(B3:G3)"

$locality.Activate() | Out-Null
$locality.Range("H6").Select() | Out-Null
